$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MBA")

# --- Populate the new "unique_id" values for rows 2-8 (column J) ---
$ws.Range("J2").Value = "VEC-017-03-175"
$ws.Range("J3").Value = "VEC-017-04-176"
$ws.Range("J4").Value = "VEC-017-04-179"
$ws.Range("J5").Value = "VEC-017-04-181"
$ws.Range("J6").Value = "VEC-017-02-188"
$ws.Range("J7").Value = "VEC-017-04-189"
$ws.Range("J8").Value = "VEC-017-04-194"

# --- Row 9 (Nirmalraj) gains a full set of profile links + unique id ---
$ws.Range("C9").Value = "https://drive.google.com/u/0/open?usp=forms_web&id=1zjRx5eM8rd-b50uPMZc4ih1lFyGdPCcB"
$ws.Range("D9").Value = "https://scholar.google.co.in/citations?user=kECK4sYAAAAJ&hl=en"
$ws.Range("E9").Value = "https://www.researchgate.net/profile/Nirmal-Amal-Raj"
$ws.Range("F9").Value = "https://orcid.org/0009-0001-6856-7412"
$ws.Range("H9").Value = "https://www.scopus.com/authid/detail.uri?authorId=59511587500"
$ws.Range("I9").Value = "https://in.linkedin.com/in/nirmal-raj-a-a8201419"
$ws.Range("J9").Value = "VEC-017-01-205"

# --- Widen column I to fit the longer links now present ---
$ws.Columns.Item(9).ColumnWidth = 57.453125

# --- Update the view: scroll position + active selection ---
$ws.Range("I10").Select()
$excel.ActiveWindow.ScrollColumn = 4
